$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.339.86"
$ws.Range("E2").Value = "  -2.75%  "
$ws.Range("D3").Value = "2.218.56"
$ws.Range("E3").Value = "  -2.05%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'107.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -9.19%  "
$ws.Range("D6").Value = "'296.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +11.80%  "
$ws.Range("D7").Value = "'0.619"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.84%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "'0.598"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.22%  "
$ws.Range("D10").Value = "'43.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.02%  "
$ws.Range("D11").Value = "'0.0912"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.05%  "
$ws.Range("D12").Value = "'54.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("E13").Value = "  -4.03%  "
$ws.Range("D14").Value = "'0.993"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.16%  "
$ws.Range("E15").Value = "  -2.57%  "
$ws.Range("D16").Value = "'15.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.04%  "
$ws.Range("D17").Value = "2.550.20"
$ws.Range("E17").Value = "  -2.24%  "
$ws.Range("D18").Value = "2.235.98"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").Value = "42.226.03"
$ws.Range("E19").Value = "  -2.87%  "
$ws.Range("D20").Value = "'7.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.73%  "
$ws.Range("E21").Value = "  -4.03%  "
$ws.Range("D22").Value = "'72.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").Value = "'3.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +20.60%  "
$ws.Range("D24").Value = "'2.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.93%  "
$ws.Range("D25").Value = "'228.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.90%  "
$ws.Range("D26").Value = "'9.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.02%  "
$ws.Range("E27").Value = "  -1.84%  "
$ws.Range("D28").Value = "'11.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.63%  "
$ws.Range("D30").Value = "'38.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.09%  "
$ws.Range("B31").Value = "WEMIXToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D31").Value = "'3.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.24%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "'173.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.47%  "
$ws.Range("D33").Value = "'20.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.88%  "
$ws.Range("D34").Value = "'0.0898"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.52%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'5.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.93%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "'5.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.24%  "
$ws.Range("D37").Value = "'4.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.53%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0381"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "'0.125"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.12%  "
$ws.Range("D40").Value = "'0.103"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.99%  "
$ws.Range("D41").Value = "'2.40"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.13%  "
$ws.Range("D42").Value = "'71.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.81%  "
$ws.Range("D43").Value = "'0.231"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.12%  "
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("D45").Value = "'12.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.01%  "
$ws.Range("E46").Value = "  -4.43%  "
$ws.Range("E47").Value = "  -6.80%  "
$ws.Range("D48").Value = "'1.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.65%  "
$ws.Range("D49").Value = "'103.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("E50").Value = "  +6.37%  "
$ws.Range("E51").Value = "  -1.76%  "
